$wb = $excel.ActiveWorkbook

# ALC!69
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 4400
$ws.Range("I69").Value = 4400
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 13200
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -12326
$ws.Range("N69").ClearContents()

# ALC!72
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 4400
$ws.Range("I72").Value = 4400
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 39600
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -35232
$ws.Range("N72").ClearContents()

# ALC!112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1426.3636
$ws.Range("I112").Value = 1050
$ws.Range("J112").Value = 1485.7894
$ws.Range("K112").Value = 3150
$ws.Range("L112").Value = 4457.3682
$ws.Range("M112").Value = -2042
$ws.Range("N112").Value = -6673.3682

# ALC!129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1543.6471
$ws.Range("I129").Value = 546
$ws.Range("J129").Value = 1676.6666
$ws.Range("K129").Value = 1638
$ws.Range("L129").Value = 5029.9998
$ws.Range("M129").Value = 3362
$ws.Range("N129").Value = -15029.9998

# ALC!132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 419168.62
$ws.Range("I132").Value = 3079.125
$ws.Range("J132").Value = 1251347.6
$ws.Range("K132").Value = 9237.375
$ws.Range("L132").Value = 3754042.8
$ws.Range("M132").Value = -6707.375
$ws.Range("N132").Value = -3759102.8

# ARM!61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1076.6129
$ws.Range("I61").Value = 1060.48
$ws.Range("J61").Value = 1143.8334
$ws.Range("K61").Value = 1060.48
$ws.Range("L61").Value = 1143.8334
$ws.Range("M61").Value = -848.48
$ws.Range("N61").Value = -1567.8334

# ARM!97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 570
$ws.Range("I97").Value = 570
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 570
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -74
$ws.Range("N97").ClearContents()

# ARM!110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 976.1429000000001
$ws.Range("I110").Value = 1006.6
$ws.Range("J110").Value = 900
$ws.Range("K110").Value = 1006.6
$ws.Range("L110").Value = 900
$ws.Range("M110").Value = 1038.4
$ws.Range("N110").Value = -4990

# ARM!136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1076.6129
$ws.Range("I136").Value = 1060.48
$ws.Range("J136").Value = 1143.8334
$ws.Range("K136").Value = 3181.44
$ws.Range("L136").Value = 3431.5002
$ws.Range("M136").Value = -631.4400000000001
$ws.Range("N136").Value = -8531.5002

# ARM!139
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 36000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 36000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 36000
$ws.Range("N139").Value = -46280

# BSM!81
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 21125
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 21125
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 21125
$ws.Range("N81").Value = -23247

# BSM!84
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H84").Value = 21125
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 21125
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 63375
$ws.Range("N84").Value = -73983

# BSM!86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2164.8696
$ws.Range("I86").Value = 1938.4445
$ws.Range("J86").Value = 2980
$ws.Range("K86").Value = 1938.4445
$ws.Range("L86").Value = 2980
$ws.Range("M86").Value = -815.4445000000001
$ws.Range("N86").Value = -5226

# BSM!89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2164.8696
$ws.Range("I89").Value = 1938.4445
$ws.Range("J89").Value = 2980
$ws.Range("K89").Value = 9692.2225
$ws.Range("L89").Value = 14900
$ws.Range("M89").Value = -4076.2225
$ws.Range("N89").Value = -26132

# CRP!16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1116.7368
$ws.Range("I16").Value = 1087.6666
$ws.Range("J16").Value = 1225.75
$ws.Range("K16").Value = 1087.6666
$ws.Range("L16").Value = 1225.75
$ws.Range("M16").Value = -800.6666
$ws.Range("N16").Value = -1799.75

# CRP!31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1994.5405
$ws.Range("I31").Value = 2022.1666
$ws.Range("J31").Value = 1000
$ws.Range("K31").Value = 2022.1666
$ws.Range("L31").Value = 1000
$ws.Range("M31").Value = -1727.1666
$ws.Range("N31").Value = -1590

# CRP!34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1994.5405
$ws.Range("I34").Value = 2022.1666
$ws.Range("J34").Value = 1000
$ws.Range("K34").Value = 2022.1666
$ws.Range("L34").Value = 1000
$ws.Range("M34").Value = -1820.1666
$ws.Range("N34").Value = -1404

# CRP!113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1116.7368
$ws.Range("I113").Value = 1087.6666
$ws.Range("J113").Value = 1225.75
$ws.Range("K113").Value = 1087.6666
$ws.Range("L113").Value = 1225.75
$ws.Range("M113").Value = 1082.3334
$ws.Range("N113").Value = -5565.75

# CUL!92
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 620.75

# CUL!119
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 1598.6666
$ws.Range("I119").Value = 1598.6666
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 4795.9998
$ws.Range("L119").Value = 0
$ws.Range("M119").Value = 42.0002000000004
$ws.Range("N119").ClearContents()

# GSM!70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4211.2
$ws.Range("I70").Value = 3931.6843
$ws.Range("J70").Value = 5728.5713
$ws.Range("K70").Value = 3931.6843
$ws.Range("L70").Value = 5728.5713
$ws.Range("M70").Value = -3661.6843
$ws.Range("N70").Value = -6268.5713

# GSM!73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4211.2
$ws.Range("I73").Value = 3931.6843
$ws.Range("J73").Value = 5728.5713
$ws.Range("K73").Value = 3931.6843
$ws.Range("L73").Value = 5728.5713
$ws.Range("M73").Value = -2995.6843
$ws.Range("N73").Value = -7600.5713

# GSM!97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1977.625
$ws.Range("I97").Value = 2121.7646
$ws.Range("J97").Value = 1627.5714
$ws.Range("K97").Value = 2121.7646
$ws.Range("L97").Value = 1627.5714
$ws.Range("M97").Value = -1625.7646
$ws.Range("N97").Value = -2619.5714

# GSM!102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1838
$ws.Range("I102").Value = 841
$ws.Range("J102").Value = 2835
$ws.Range("K102").Value = 841
$ws.Range("L102").Value = 2835
$ws.Range("M102").Value = 781
$ws.Range("N102").Value = -6079

# LTW!68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2726.8
$ws.Range("I68").Value = 2643.1428
$ws.Range("J68").Value = 2800
$ws.Range("K68").Value = 2643.1428
$ws.Range("L68").Value = 2800
$ws.Range("M68").Value = -1894.1428
$ws.Range("N68").Value = -4298

# LTW!71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2726.8
$ws.Range("I71").Value = 2643.1428
$ws.Range("J71").Value = 2800
$ws.Range("K71").Value = 13215.714
$ws.Range("L71").Value = 14000
$ws.Range("M71").Value = -9471.714
$ws.Range("N71").Value = -21488

# LTW!82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2741.6667
$ws.Range("I82").Value = 2500
$ws.Range("J82").Value = 3080
$ws.Range("K82").Value = 2500
$ws.Range("L82").Value = 3080
$ws.Range("M82").Value = -2139
$ws.Range("N82").Value = -3802

# LTW!85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2741.6667
$ws.Range("I85").Value = 2500
$ws.Range("J85").Value = 3080
$ws.Range("K85").Value = 2500
$ws.Range("L85").Value = 3080
$ws.Range("M85").Value = -1252
$ws.Range("N85").Value = -5576

# LTW!136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2500.9395
$ws.Range("I136").Value = 1387.7931
$ws.Range("J136").Value = 10571.25
$ws.Range("K136").Value = 4163.379300000001
$ws.Range("L136").Value = 31713.75
$ws.Range("M136").Value = -1613.379300000001
$ws.Range("N136").Value = -36813.75

# WVR!132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 975.7963
$ws.Range("I132").Value = 767.5
$ws.Range("J132").Value = 2173.5
$ws.Range("K132").Value = 2302.5
$ws.Range("L132").Value = 6520.5
$ws.Range("M132").Value = 227.5
$ws.Range("N132").Value = -11580.5
